$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the ERSST_V5 Temperature value and the two anomaly columns that
# derive from it (Uk37 Temperature anomaly_Muller98/BAYSPLINE - ERSST).
$ws.Range("P2").Value = 17.76
$ws.Range("Q2").Value = 2.811264970567493
$ws.Range("R2").Value = 2.568582730567492

# Remove the "MgCa Temperature anomaly_Original - Coretop" column entirely
# (header in Z1, value in Z2), shrinking the used range to A1:Y2.
$ws.Range("Z1:Z2").Clear()
